$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Double-space all three existing paragraphs (w:spacing w:line="480"
#    w:lineRule="auto" == wdLineSpaceDouble).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $p.Range.ParagraphFormat.LineSpacingRule = 2
}

# ---------------------------------------------------------------------
# 2) Append the new literature-review sentences to the end of the third
#    paragraph, right before the "_GoBack" bookmark, preserving the
#    run-level formatting (yellow highlight / italics) from the diff.
# ---------------------------------------------------------------------
function Append-Run([string]$text, [bool]$highlight, [bool]$italic) {
    $bm = $d.Bookmarks("_GoBack")
    $startPos = $bm.Start
    $r = $d.Range($startPos, $startPos)
    $r.InsertBefore($text)

    $bm2 = $d.Bookmarks("_GoBack")
    $endPos = $bm2.Start
    $rng = $d.Range($startPos, $endPos)

    if ($highlight) {
        $rng.Font.HighlightColorIndex = 7
    }
    if ($italic) {
        $rng.Italic = 1
    }
}

Append-Run " " $true $false
Append-Run " " $false $false
Append-Run "Other terms such as vulnerability, accessibility, and connectivity are all terms that factor into resiliency as a whole. However, this paper will seek to discuss only how these terms relate to overall resiliency. Accessibility refers to the ease of which links or routes are able to connect to nodes in a travel model (" $false $false
Append-Run "Cantillo" $true $false
Append-Run " et al.)." $true $false
Append-Run " Similarly, connectivity is defined as how extensively a network is " $false $false
Append-Run "inter" $false $true
Append-Run "connected. " $false $false
Append-Run "Cardinale" $true $false
Append-Run " suggests that connectivity often determines how susceptible a network is to disruption." $false $false
